$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D to text format so numeric-looking strings
# (e.g. "585.49") are stored as text instead of being parsed as numbers,
# matching the original inline-string cell content.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "67.201.01"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "2.483.77"
$ws.Range("E3").Value = "  +0.53%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "585.49"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").Value = "172.68"
$ws.Range("E6").Value = "  +3.05%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +0.30%  "
$ws.Range("D9").Value = "2.483.25"
$ws.Range("E9").Value = "  +0.48%  "
$ws.Range("E10").Value = "  +2.99%  "
$ws.Range("E11").Value = "  +1.14%  "
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("D13").Value = "0.334"
$ws.Range("E13").Value = "  +0.14%  "
$ws.Range("D14").Value = "2.931.94"
$ws.Range("D15").Value = "25.57"
$ws.Range("E15").Value = "  +0.35%  "
$ws.Range("D16").Value = "67.041.78"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("E17").Value = "  +1.00%  "
$ws.Range("D18").Value = "2.493.60"
$ws.Range("E18").Value = "  -0.48%  "
$ws.Range("D19").Value = "7.59"
$ws.Range("E19").Value = "  -0.87%  "
$ws.Range("D20").Value = "10.98"
$ws.Range("E20").Value = "  -3.35%  "
$ws.Range("D21").Value = "350.60"
$ws.Range("E21").Value = "  -1.50%  "
$ws.Range("E22").Value = "  -0.75%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("E24").Value = "  -0.79%  "
$ws.Range("E25").Value = "  +0.28%  "
$ws.Range("E26").Value = "  +3.02%  "
$ws.Range("E27").Value = "  +1.83%  "
$ws.Range("D28").Value = "2.609.02"
$ws.Range("E28").Value = "  +0.56%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.48%  "
$ws.Range("D30").Value = "0.0₃0913"
$ws.Range("E30").Value = "  +1.36%  "
$ws.Range("D31").Value = "507.80"
$ws.Range("E31").Value = "  -1.07%  "
$ws.Range("D32").Value = "7.71"
$ws.Range("E32").Value = "  -1.10%  "
$ws.Range("E33").Value = "  +1.33%  "
$ws.Range("E34").Value = "  -0.72%  "
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").Value = "162.41"
$ws.Range("E36").Value = "  +2.33%  "
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("E38").Value = "  +0.59%  "
$ws.Range("E39").Value = "  -1.70%  "
$ws.Range("D40").Value = "1.34"
$ws.Range("E40").Value = "  -0.37%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("E42").Value = "  +1.41%  "
$ws.Range("E43").Value = "  +1.14%  "
$ws.Range("D44").Value = "4.83"
$ws.Range("E44").Value = "  +0.86%  "
$ws.Range("D45").Value = "2.39"
$ws.Range("E45").Value = "  +3.13%  "
$ws.Range("D46").Value = "143.49"
$ws.Range("E46").Value = "  +1.39%  "
$ws.Range("E47").Value = "  +3.21%  "
$ws.Range("D48").Value = "3.49"
$ws.Range("E48").Value = "  +0.32%  "
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("E50").Value = "  +0.41%  "
$ws.Range("E51").Value = "  -0.90%  "

# Restore default (General) formatting/style on column D so the
# cells keep their original style index while remaining text cells.
$ws.Range("D2:D51").ClearFormats()
